$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'FAPs'
$ws.Cells.Item(2, 2).Value = 'Sectm1a'
$ws.Cells.Item(2, 3).Value = 'Cd7'
$ws.Cells.Item(2, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.175109
$ws.Cells.Item(2, 8).Value = 0.525327
$ws.Cells.Item(2, 9).Value = 0.641283184058224
$ws.Cells.Item(2, 10).Value = 0.641283184058224
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.5923286666666666
$ws.Cells.Item(2, 14).Value = 1.776986
$ws.Cells.Item(2, 15).Value = 0.2489148585411457
$ws.Cells.Item(2, 16).Value = 0.2489148585411457
$ws.Cells.Item(2, 17).Value = 0.1037220804913333
$ws.Cells.Item(2, 18).Value = 0.933498724422
$ws.Cells.Item(2, 19).Value = 0.1596249130446683
$ws.Cells.Item(2, 20).Value = 0.1596249130446683

# Row 3
$ws.Cells.Item(3, 1).Value = 'FAPs'
$ws.Cells.Item(3, 2).Value = 'Sectm1a'
$ws.Cells.Item(3, 3).Value = 'Cd7'
$ws.Cells.Item(3, 4).Value = 'Neutrophils'
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.175109
$ws.Cells.Item(3, 8).Value = 0.525327
$ws.Cells.Item(3, 9).Value = 0.641283184058224
$ws.Cells.Item(3, 10).Value = 0.641283184058224
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.387020666666667
$ws.Cells.Item(3, 14).Value = 4.161062
$ws.Cells.Item(3, 15).Value = 0.5828690598074137
$ws.Cells.Item(3, 16).Value = 0.5828690598074138
$ws.Cells.Item(3, 17).Value = 0.2428798019193333
$ws.Cells.Item(3, 18).Value = 2.185918217274
$ws.Cells.Item(3, 19).Value = 0.3737841265623216
$ws.Cells.Item(3, 20).Value = 0.3737841265623217

# Row 4
$ws.Cells.Item(4, 1).Value = 'FAPs'
$ws.Cells.Item(4, 2).Value = 'Sectm1a'
$ws.Cells.Item(4, 3).Value = 'Cd7'
$ws.Cells.Item(4, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.175109
$ws.Cells.Item(4, 8).Value = 0.525327
$ws.Cells.Item(4, 9).Value = 0.641283184058224
$ws.Cells.Item(4, 10).Value = 0.641283184058224
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.4002943333333334
$ws.Cells.Item(4, 14).Value = 1.200883
$ws.Cells.Item(4, 15).Value = 0.1682160816514405
$ws.Cells.Item(4, 16).Value = 0.1682160816514405
$ws.Cells.Item(4, 17).Value = 0.07009514041566667
$ws.Cells.Item(4, 18).Value = 0.630856263741
$ws.Cells.Item(4, 19).Value = 0.107874144451234
$ws.Cells.Item(4, 20).Value = 0.107874144451234

# Row 5
$ws.Cells.Item(5, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(5, 2).Value = 'Sectm1a'
$ws.Cells.Item(5, 3).Value = 'Cd7'
$ws.Cells.Item(5, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.06976733333333333
$ws.Cells.Item(5, 8).Value = 0.209302
$ws.Cells.Item(5, 9).Value = 0.2555015314075889
$ws.Cells.Item(5, 10).Value = 0.2555015314075888
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.5923286666666666
$ws.Cells.Item(5, 14).Value = 1.776986
$ws.Cells.Item(5, 15).Value = 0.2489148585411457
$ws.Cells.Item(5, 16).Value = 0.2489148585411457
$ws.Cells.Item(5, 17).Value = 0.04132519153022222
$ws.Cells.Item(5, 18).Value = 0.371926723772
$ws.Cells.Item(5, 19).Value = 0.06359812754736607
$ws.Cells.Item(5, 20).Value = 0.06359812754736606

# Row 6
$ws.Cells.Item(6, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(6, 2).Value = 'Sectm1a'
$ws.Cells.Item(6, 3).Value = 'Cd7'
$ws.Cells.Item(6, 4).Value = 'Neutrophils'
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.06976733333333333
$ws.Cells.Item(6, 8).Value = 0.209302
$ws.Cells.Item(6, 9).Value = 0.2555015314075889
$ws.Cells.Item(6, 10).Value = 0.2555015314075888
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.387020666666667
$ws.Cells.Item(6, 14).Value = 4.161062
$ws.Cells.Item(6, 15).Value = 0.5828690598074137
$ws.Cells.Item(6, 16).Value = 0.5828690598074138
$ws.Cells.Item(6, 17).Value = 0.09676873319155556
$ws.Cells.Item(6, 18).Value = 0.870918598724
$ws.Cells.Item(6, 19).Value = 0.1489239373908957
$ws.Cells.Item(6, 20).Value = 0.1489239373908957

# Row 7
$ws.Cells.Item(7, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(7, 2).Value = 'Sectm1a'
$ws.Cells.Item(7, 3).Value = 'Cd7'
$ws.Cells.Item(7, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.06976733333333333
$ws.Cells.Item(7, 8).Value = 0.209302
$ws.Cells.Item(7, 9).Value = 0.2555015314075889
$ws.Cells.Item(7, 10).Value = 0.2555015314075888
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.4002943333333334
$ws.Cells.Item(7, 14).Value = 1.200883
$ws.Cells.Item(7, 15).Value = 0.1682160816514405
$ws.Cells.Item(7, 16).Value = 0.1682160816514405
$ws.Cells.Item(7, 17).Value = 0.02792746818511111
$ws.Cells.Item(7, 18).Value = 0.251347213666
$ws.Cells.Item(7, 19).Value = 0.04297946646932706
$ws.Cells.Item(7, 20).Value = 0.04297946646932706

# Row 8
$ws.Cells.Item(8, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(8, 2).Value = 'Sectm1a'
$ws.Cells.Item(8, 3).Value = 'Cd7'
$ws.Cells.Item(8, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.028184
$ws.Cells.Item(8, 8).Value = 0.084552
$ws.Cells.Item(8, 9).Value = 0.1032152845341872
$ws.Cells.Item(8, 10).Value = 0.1032152845341872
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.5923286666666666
$ws.Cells.Item(8, 14).Value = 1.776986
$ws.Cells.Item(8, 15).Value = 0.2489148585411457
$ws.Cells.Item(8, 16).Value = 0.2489148585411457
$ws.Cells.Item(8, 17).Value = 0.01669419114133333
$ws.Cells.Item(8, 18).Value = 0.150247720272
$ws.Cells.Item(8, 19).Value = 0.02569181794911131
$ws.Cells.Item(8, 20).Value = 0.02569181794911131

# Row 9
$ws.Cells.Item(9, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(9, 2).Value = 'Sectm1a'
$ws.Cells.Item(9, 3).Value = 'Cd7'
$ws.Cells.Item(9, 4).Value = 'Neutrophils'
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.028184
$ws.Cells.Item(9, 8).Value = 0.084552
$ws.Cells.Item(9, 9).Value = 0.1032152845341872
$ws.Cells.Item(9, 10).Value = 0.1032152845341872
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.387020666666667
$ws.Cells.Item(9, 14).Value = 4.161062
$ws.Cells.Item(9, 15).Value = 0.5828690598074137
$ws.Cells.Item(9, 16).Value = 0.5828690598074138
$ws.Cells.Item(9, 17).Value = 0.03909179046933334
$ws.Cells.Item(9, 18).Value = 0.351826114224
$ws.Cells.Item(9, 19).Value = 0.06016099585419638
$ws.Cells.Item(9, 20).Value = 0.06016099585419639

# Row 10
$ws.Cells.Item(10, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(10, 2).Value = 'Sectm1a'
$ws.Cells.Item(10, 3).Value = 'Cd7'
$ws.Cells.Item(10, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.028184
$ws.Cells.Item(10, 8).Value = 0.084552
$ws.Cells.Item(10, 9).Value = 0.1032152845341872
$ws.Cells.Item(10, 10).Value = 0.1032152845341872
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.4002943333333334
$ws.Cells.Item(10, 14).Value = 1.200883
$ws.Cells.Item(10, 15).Value = 0.1682160816514405
$ws.Cells.Item(10, 16).Value = 0.1682160816514405
$ws.Cells.Item(10, 17).Value = 0.01128189549066667
$ws.Cells.Item(10, 18).Value = 0.101537059416
$ws.Cells.Item(10, 19).Value = 0.0173624707308795
$ws.Cells.Item(10, 20).Value = 0.0173624707308795

# Remove now-obsolete rows 11-13 (old Resolving-Mac block superseded by rows 8-10)
$ws.Range("A11:T13").EntireRow.Delete()
